# Generate Report for Handback
# Updates the handback-status report with refreshed timestamps and priority.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for rows 2 and 4
$wsOverview.Range("G2").Value = "2016-08-17 10:15:15"
$wsOverview.Range("G4").Value = "2016-08-17 10:15:15"

# zh-cn sheet: Priority column (E) rows 2 and 4: ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime column (H) rows 2 and 4
$wsZhCn.Range("H2").Value = "2016-08-17 10:15:00"
$wsZhCn.Range("H4").Value = "2016-08-17 10:15:00"

# zh-cn sheet: Correspond Handback DateTime column (K) rows 2 and 4
$wsZhCn.Range("K2").Value = "2016-08-17 10:15:31"
$wsZhCn.Range("K4").Value = "2016-08-17 10:15:31"

# de-de sheet: Priority column (E) rows 2 and 4: ht -> mt
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# de-de sheet: Correspond Handoff Datetime column (H) rows 2 and 4
$wsDeDe.Range("H2").Value = "2016-08-17 10:15:15"
$wsDeDe.Range("H4").Value = "2016-08-17 10:15:15"

# de-de sheet: Correspond Handback DateTime column (K) rows 2 and 4
$wsDeDe.Range("K2").Value = "2016-08-17 10:15:38"
$wsDeDe.Range("K4").Value = "2016-08-17 10:15:38"
